$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the table by copying row 186 formatting down to new row 187
$ws.Range("A186:DX186").Copy()
$ws.Range("A187:DX187").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# Populate new row 187 with data
$ws.Range("A187").Value = 44081
$ws.Range("B187").Value = 671848
$ws.Range("C187").Value = 2714
$ws.Range("D187").Value = 89341
$ws.Range("E187").Value = 65113
$ws.Range("F187").Value = 227976
$ws.Range("G187").Value = 26015
$ws.Range("H187").Value = 4417
$ws.Range("I187").Value = 3466
$ws.Range("J187").Value = 6897
$ws.Range("K187").Value = 6409
$ws.Range("L187").Value = 13010
$ws.Range("M187").Value = 3789
$ws.Range("N187").Value = 21139
$ws.Range("O187").Value = 26038
$ws.Range("P187").Value = 5934
$ws.Range("Q187").Value = 6260
$ws.Range("R187").Value = 13143
$ws.Range("S187").Value = 10227
$ws.Range("T187").Value = 15271
$ws.Range("U187").Value = 12880
$ws.Range("V187").Value = 3247
$ws.Range("W187").Value = 1617
$ws.Range("X187").Value = 7078
$ws.Range("Y187").Value = 21912
$ws.Range("Z187").Value = 12543
$ws.Range("AA187").Value = 8491
$ws.Range("AB187").Value = 50420
$ws.Range("AC187").Value = 1294
$ws.Range("AD187").Value = 372
$ws.Range("AE187").Value = 448
$ws.Range("AF187").Value = 454
$ws.Range("AG187").Value = 277
$ws.Range("AH187").Value = 202
$ws.Range("AI187").Value = 352
$ws.Range("AJ187").Value = 1985
$ws.Range("AK187").Value = 3762
$ws.Range("AL187").Value = 36684
$ws.Range("AM187").Value = 7625
$ws.Range("AN187").Value = 2432
$ws.Range("AO187").Value = 39512
$ws.Range("AP187").Value = 1009
$ws.Range("AQ187").Value = 20946
$ws.Range("AR187").Value = 1475
$ws.Range("AS187").Value = 8916
$ws.Range("AT187").Value = 1564
$ws.Range("AU187").Value = 1582
$ws.Range("AV187").Value = 5498
$ws.Range("AW187").Value = 1699
$ws.Range("AX187").Value = 951
$ws.Range("AY187").Value = 2481
$ws.Range("AZ187").Value = 2647
$ws.Range("BA187").Value = 51727
$ws.Range("BB187").Value = 12905
$ws.Range("BC187").Value = 3554
$ws.Range("BD187").Value = 8226
$ws.Range("BE187").Value = 4715
$ws.Range("BF187").Value = 280
$ws.Range("BG187").Value = 1418
$ws.Range("BH187").Value = 2626
$ws.Range("BI187").Value = 734
$ws.Range("BJ187").Value = 2067
$ws.Range("BK187").Value = 8795
$ws.Range("BL187").Value = 8854
$ws.Range("BM187").Value = 9194
$ws.Range("BN187").Value = 13969
$ws.Range("BO187").Value = 1893
$ws.Range("BP187").Value = 835
$ws.Range("BQ187").Value = 9471
$ws.Range("BR187").Value = 8139
$ws.Range("BS187").Value = 9537
$ws.Range("BT187").Value = 1822
$ws.Range("BU187").Value = 1692
$ws.Range("BV187").Value = 3827
$ws.Range("BW187").Value = 3829
$ws.Range("BX187").Value = 1165
$ws.Range("BY187").Value = 4987
$ws.Range("BZ187").Value = 2739
$ws.Range("CA187").Value = 1489
$ws.Range("CB187").Value = 800
$ws.Range("CC187").Value = 2381
$ws.Range("CD187").Value = 2036
$ws.Range("CE187").Value = 1502
$ws.Range("CF187").Value = 1116
$ws.Range("CG187").Value = 5597
$ws.Range("CH187").Value = 1638
$ws.Range("CI187").Value = 1251
$ws.Range("CJ187").Value = 1426
$ws.Range("CK187").Value = 1810
$ws.Range("CL187").Value = 1690
$ws.Range("CM187").Value = 2020
$ws.Range("CN187").Value = 1299
$ws.Range("CO187").Value = 1113
$ws.Range("CP187").Value = 1131
$ws.Range("CQ187").Value = 665
$ws.Range("CR187").Value = 3121
$ws.Range("CS187").Value = 1169
$ws.Range("CT187").Value = 827
$ws.Range("CU187").Value = 811
$ws.Range("CV187").Value = 1478
$ws.Range("CW187").Value = 1344
$ws.Range("CX187").Value = 691
$ws.Range("CY187").Value = 783
$ws.Range("CZ187").Value = 1047
$ws.Range("DA187").Value = 1313
$ws.Range("DB187").Value = 1140
$ws.Range("DC187").Value = 1263
$ws.Range("DD187").Value = 974
$ws.Range("DE187").Value = 319
$ws.Range("DF187").Value = 344
$ws.Range("DG187").Value = 734
$ws.Range("DH187").Value = 652
$ws.Range("DI187").Value = 429
$ws.Range("DJ187").Value = 534
$ws.Range("DK187").Value = 352
$ws.Range("DL187").Value = 629
$ws.Range("DM187").Value = 720
$ws.Range("DN187").Value = 517
$ws.Range("DO187").Value = 481
$ws.Range("DP187").Value = 372
$ws.Range("DQ187").Value = 517
$ws.Range("DR187").Value = 123215
$ws.Range("DS187").Value = 285130
$ws.Range("DT187").Value = 12434
$ws.Range("DU187").Value = 123013
$ws.Range("DV187").Value = 75889
$ws.Range("DW187").Value = 34841
$ws.Range("DX187").Value = 10421

# Apply individual cell corrections throughout the sheet
$ws.Range("BS11").Value = "NaN"
$ws.Range("DW11").Value = 2
$ws.Range("BM13").Value = 1
$ws.Range("BM14").Value = 1
$ws.Range("BS14").Value = 1
$ws.Range("BM15").Value = 2
$ws.Range("BM16").Value = 4
$ws.Range("BM17").Value = 5
$ws.Range("BM18").Value = 5
$ws.Range("BM19").Value = 6
$ws.Range("BM20").Value = 6
$ws.Range("BM21").Value = 7
$ws.Range("AK33").Value = "NaN"
$ws.Range("AK35").Value = 1
$ws.Range("DC35").Value = "NaN"
$ws.Range("DC37").Value = 1
$ws.Range("BM64").Value = 88
$ws.Range("CF93").Value = "NaN"
$ws.Range("CF132").Value = "NaN"
$ws.Range("BM89").Value = 326
$ws.Range("BM90").Value = 343
$ws.Range("BM91").Value = 354
$ws.Range("BM92").Value = 376
$ws.Range("BM93").Value = 384
$ws.Range("BM94").Value = 401
$ws.Range("BM95").Value = 420
$ws.Range("BM96").Value = 449
$ws.Range("BM97").Value = 469
$ws.Range("BM98").Value = 473
$ws.Range("BM99").Value = 496
$ws.Range("BM100").Value = 519
$ws.Range("BM101").Value = 556
$ws.Range("BM102").Value = 569
$ws.Range("BM103").Value = 580
$ws.Range("BM104").Value = 603
$ws.Range("BM105").Value = 615
$ws.Range("BM106").Value = 646
$ws.Range("BM107").Value = 675
$ws.Range("BM108").Value = 687
$ws.Range("BM109").Value = 700
$ws.Range("BM110").Value = 722
$ws.Range("BM111").Value = 772
$ws.Range("BM112").Value = 797
$ws.Range("BM113").Value = 826
$ws.Range("BM114").Value = 870
$ws.Range("BM115").Value = 908
$ws.Range("BM116").Value = 941
$ws.Range("BM117").Value = 953
$ws.Range("BM118").Value = 1021
$ws.Range("BM119").Value = 1071
$ws.Range("BM120").Value = 1127
$ws.Range("BM121").Value = 1156
$ws.Range("BM122").Value = 1211
$ws.Range("BM123").Value = 1230
$ws.Range("BM124").Value = 1235
$ws.Range("BM125").Value = 1244
$ws.Range("BM126").Value = 1299
$ws.Range("BM127").Value = 1360
$ws.Range("BM128").Value = 1518
$ws.Range("BM129").Value = 1552
$ws.Range("BM130").Value = 1607
$ws.Range("BM131").Value = 1678
$ws.Range("BM132").Value = 1739
$ws.Range("BM133").Value = 1809
$ws.Range("BM134").Value = 1944
$ws.Range("BM135").Value = 2063
$ws.Range("BM136").Value = 2174
$ws.Range("BM137").Value = 2226
$ws.Range("BM138").Value = 2300
$ws.Range("BM139").Value = 2344
$ws.Range("BM140").Value = 2449
$ws.Range("BM141").Value = 2523
$ws.Range("BM142").Value = 2647
$ws.Range("BM143").Value = 2680
$ws.Range("BM144").Value = 2759
$ws.Range("BM145").Value = 2795
$ws.Range("BM146").Value = 2892
$ws.Range("BM147").Value = 2994
$ws.Range("BM148").Value = 3214
$ws.Range("BM149").Value = 3284
$ws.Range("BM150").Value = 3499
$ws.Range("BM151").Value = 3716
$ws.Range("BM152").Value = 3832
$ws.Range("BM153").Value = 3895
$ws.Range("BM154").Value = 4052
$ws.Range("BM155").Value = 4143
$ws.Range("BM156").Value = 4240
$ws.Range("BM157").Value = 4343
$ws.Range("BM158").Value = 4395
$ws.Range("BM159").Value = 4470
$ws.Range("BM160").Value = 4643
$ws.Range("BM161").Value = 4938
$ws.Range("BM162").Value = 5167
$ws.Range("BM163").Value = 5275
$ws.Range("BM164").Value = 5503
$ws.Range("BM165").Value = 5815
$ws.Range("BM166").Value = 5929
$ws.Range("BM167").Value = 6240
$ws.Range("BM168").Value = 6445
$ws.Range("BM169").Value = 6944
$ws.Range("BM170").Value = 7054
$ws.Range("BM171").Value = 7124
$ws.Range("BM172").Value = 7221
$ws.Range("BM173").Value = 7474
$ws.Range("BM174").Value = 7686
$ws.Range("BM175").Value = 7736
$ws.Range("BM176").Value = 7909
$ws.Range("BM177").Value = 8019
$ws.Range("BM178").Value = 8189
$ws.Range("BM179").Value = 8364

# Restore selection to match final state
$ws.Range("A186").Select()
